# GPIO.xlsx update — "presentatie en exel pins"
# Adds the four new "Motor .. euro uit" labels in column D for pins 26-33
# (rows 30-37), matches their fill style to the existing highlighted D
# cells, widens column D to fit the new text, and updates the sheet's
# active selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new motor-off labels (D30:D37) ------------------------------------
$ws.Range("D30").Value = "Motor 5 euro uit"
$ws.Range("D31").Value = "Motor 5 euro uit"
$ws.Range("D32").Value = "Motor 10 euro uit"
$ws.Range("D33").Value = "Motor 10 euro uit"
$ws.Range("D34").Value = "Motor 20 euro uit"
$ws.Range("D35").Value = "Motor 20 euro uit"
$ws.Range("D36").Value = "Motor 50 euro uit"
$ws.Range("D37").Value = "Motor 50 euro uit"

# Match the highlighted fill style already used on the other D-column
# annotations (e.g. D17, D19:D27) by copying their formatting over.
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D30:D37").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- column width so the new labels fit ---------------------------------
$ws.Columns("D").ColumnWidth = 16.3

# --- view state: selection + scroll position -----------------------------
$ws.Range("A16").Select() | Out-Null
$ws.Range("I34").Select() | Out-Null
